$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing base/projected dates for "Proyecto B" (row 3) -
# these cells already carry the date number format, so a plain value
# write is enough.
$ws.Range("B3").Value = 45992
$ws.Range("C3").Value = 45996

# Fill in the missing base/projected dates for "Proyecto C" (row 4) and
# give them the same date formatting used by the other rows before
# writing the values so the existing style gets reused.
$ws.Range("B4:C4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B4").Value = 45992
$ws.Range("C4").Value = 45996

# Move the active selection as the user last left it
$ws.Range("B8").Select()
